$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.8068052433981242, 0.5475239231397175, 1.5059970179533422, 2.9572477800483057)"
$ws.Range("C2").Value = "JSU(-1.0932392054973943, 1.0792698804906662, 2.700763781474773, 3.816690434003931)"
$ws.Range("D2").Value = "NIG(0.8553989732161207, 0.6483924039590527, 1.3055593656737134, 3.2062767616683585)"
$ws.Range("E2").Value = "JSU(-0.9448186720902911, 1.2507624641223174, 4.328167488086834, 6.092634829257609)"
